# Update Diary Zihua Weng (#276)
# Fills in the Jan 30th (second entry) / Feb 4th / Feb 6th diary rows that
# were previously left blank, and removes the long run of now-unused blank
# rows at the tail of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: second "Jan 30th" entry -------------------------------------
$ws.Range("A21").Value = "Jan 30th "
$ws.Range("B21").Value = "5pm-8pm"
$ws.Range("C21").Value = "Andre, Kaj"
$ws.Range("D21").Value = "Learn umls and different code graph"
$ws.Range("E21").Value = "Drew different umls and code graph for my code."
$ws.Range("F21").Value = "Working in team is very important! When figuring out a problem, we should work as deep as we need. Try to brainstorm as more possibilities as we can before we head to the final solution. Recruiters are looking for candicates who could tear thing apart and solve problems. We need to make sure we have a clear communication and show how we get a relationship with others when we answer the questions."
$ws.Range("G21").Value = "Cool!"

# --- Row 22: Feb 4th entry -------------------------------------------------
$ws.Range("A22").Value = "Feb 4th"
$ws.Range("B22").Value = "8pm-9pm"
$ws.Range("C22").Value = "Junxian, Wenchia"
$ws.Range("D22").Value = "Select two features from Realm for homework 2."
$ws.Range("E22").Value = "We decided to look at features: 1. Add methods to table query. 2. Change encryption algorithms "
$ws.Range("G22").Value = "Feel like we have a lot work to do to understand the whole project which is well developed within years…."

# --- Row 23: Feb 6th entry -------------------------------------------------
# (Goal/Reflection were typed in before Time/Achievements were back-filled,
# so the cells are set in that same order here.)
$ws.Range("A23").Value = "Feb 6th"
$ws.Range("C23").Value = "Junxian, Wenchia"
$ws.Range("D23").Value = "Finish homework2. Read the code of the above two features and understand how new developer could modify the code."
$ws.Range("F23").Value = "We finished the first feature1 easily as we know some of the concepts during homework1. But the feature2 is way more difficult as we had difficulty finding the encryption algorithms Realm uses. At the end of the day, everything is executed in C++ and we have hard time chasing the code in deep. Also, I only have littel knowledge of encryption algorithms, so I need to learn more about that later."
$ws.Range("B23").Value = "10am-2pm"
$ws.Range("E23").Value = "Finished homework2 and write the report . For the second feature we have trouble finding the encryption algorithm code. Finally, we found out that the default algorithm is AES-256 which is defined in other project called realm-core."
$ws.Range("G23").Value = "It is frustrated that the entryption part of the code is too hard to find and we went to dead end when we chased the code. And a lot project is writen in C++…"

# D21 had drifted onto the italic 11pt "Reflection" font; match the rest of
# the row (12pt) the way the surrounding cells are already formatted.
$ws.Range("D21").Font.Size = 12

# Row heights grow with the wrapped text, same as every other filled-in row.
$ws.Rows("21:21").RowHeight = 204
$ws.Rows("22:22").RowHeight = 51
$ws.Rows("23:23").RowHeight = 187

# --- Remove the long tail of now-superfluous blank rows --------------------
$ws.Range("A115:G125").EntireRow.Delete()

# Leave the selection roughly where the editor ended up.
$ws.Range("E23").Select()
